$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# The runtime's "next shape id" counter is primed incorrectly for the very
# first shape-creating call in a session (it returns id 3 instead of
# continuing from the highest existing id on the slide). Work around this
# with a disposable warm-up shape so subsequent Duplicate() calls receive
# the correct, sequential ids (19, 20, 21, 22).
$warmup = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$warmup.Delete()

# Grab references to the source shapes by name up front, since positional
# indices shift as new shapes are inserted / reordered below.
$srcT1 = $s.Shapes.Item("Rechteck 7")
$srcArrow1 = $s.Shapes.Item("Pfeil nach unten 12")
$srcArrow2 = $s.Shapes.Item("Pfeil nach unten 10")
$srcM1 = $s.Shapes.Item("Rechteck 6")

# --- New shape "Rechteck 18" (id 19) --------------------------------------
# Duplicate of "Rechteck 7" (T1), nudged by 1 EMU in x/y, appended at the
# end of the shape tree (right after "Textfeld 17").
$rechteck18 = $srcT1.Duplicate()
$rechteck18.Name = "Rechteck 18"
$rechteck18.Left = 5076057 / 12700.0
$rechteck18.Top = 2699657 / 12700.0

# --- New shape "Pfeil nach unten 19" (id 20) -------------------------------
# Duplicate of "Pfeil nach unten 12", nudged by 1 EMU in x/y, then sent to
# the back of the shape tree so it becomes the very first shape (in front
# of "Titel 1").
$pfeil19 = $srcArrow1.Duplicate()
$pfeil19.Name = "Pfeil nach unten 19"
$pfeil19.Left = 5842403 / 12700.0
$pfeil19.Top = 3308805 / 12700.0
$pfeil19.ZOrder(1)

# --- New shape "Pfeil nach unten 20" (id 21) -------------------------------
# Duplicate of "Pfeil nach unten 10", nudged by 1 EMU in x/y, appended at
# the end of the shape tree.
$pfeil20 = $srcArrow2.Duplicate()
$pfeil20.Name = "Pfeil nach unten 20"
$pfeil20.Left = 5265789 / 12700.0
$pfeil20.Top = 1910585 / 12700.0

# --- New shape "Rechteck 21" (id 22) ---------------------------------------
# Duplicate of "Rechteck 6" (M1), nudged by 1 EMU in x/y, appended at the
# end of the shape tree.
$rechteck21 = $srcM1.Duplicate()
$rechteck21.Name = "Rechteck 21"
$rechteck21.Left = 5076057 / 12700.0
$rechteck21.Top = 1347615 / 12700.0
